# New layout of Event datastructure in database: add a "data type" column
# to the Users/Events/Messages field tables, split "Time" into
# "StartTime"/"EndTime", and tidy a couple of values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "data type" column for the Users table (column F)
# ---------------------------------------------------------------------
$ws.Range("F7").Value  = "String"
$ws.Range("F8").Value  = "Map<String, Boolean>"
$ws.Range("F9").Value  = "String"
$ws.Range("F10").Value = "String"

# ---------------------------------------------------------------------
# New "data type" column for the Events table (column K)
# ---------------------------------------------------------------------
$ws.Range("K7").Value  = "String"
$ws.Range("K8").Value  = "Long"
$ws.Range("K9").Value  = "Long"
$ws.Range("K10").Value = "Map<String, Boolean>"
$ws.Range("K11").Value = "Map<String, Boolean>"
$ws.Range("K12").Value = "String"
$ws.Range("K13").Value = "String"
$ws.Range("K14").Value = "String"
$ws.Range("K15").Value = "String"
$ws.Range("K16").Value = "String"
$ws.Range("K17").Value = "String"
$ws.Range("K18").Value = "String"
$ws.Range("K19").Value = "String"
$ws.Range("K20").Value = "Map<String, Boolean>"

# ---------------------------------------------------------------------
# "Time" is replaced by separate "StartTime" / "EndTime" fields, which
# pushes OwnerOfTheEvent / Title / WhoReported down by one row.
# ---------------------------------------------------------------------
$ws.Range("I16").Value = "StartTime"
# J16 already holds 0.875 with the existing time number format - keep it.

$ws.Range("I17").Value = "EndTime"
$ws.Range("J17").ClearFormats()
$ws.Range("J17").Value = 1
$ws.Range("J17").NumberFormat = "h:mm"

$ws.Range("I18").Value = "OwnerOfTheEvent"
$ws.Range("J18").Value = "asfasdafsfasdasd (EUID)"

$ws.Range("I19").Value = "Title"
$ws.Range("J19").Value = "Nameoftheevent"

$ws.Range("I20").Value = "WhoReported"
$ws.Range("J20").Value = "ytryt,asd,asfgfds,asdasd"
$ws.Range("J20").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Column widths (approximate Mac-Excel re-save widening)
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth  = 27.666666666666668
$ws.Columns.Item(8).ColumnWidth  = 18.498697916666668
$ws.Columns.Item(10).ColumnWidth = 32.498697916666664
$ws.Columns.Item(13).ColumnWidth = 16.666666666666668
$ws.Columns.Item(14).ColumnWidth = 15.830729166666666
$ws.Columns.Item(17).ColumnWidth = 14.666666666666666

# ---------------------------------------------------------------------
# Selection, as saved in the new workbook
# ---------------------------------------------------------------------
$ws.Range("F6").Select()
